# Applies the "Updated cryptos list" data refresh to the active worksheet.
# Column D (Price) holds values formatted as text (e.g. "61.154.32"); some of those
# look like plain numbers (e.g. "549.51"), so a leading apostrophe is used -- exactly as
# typing '549.51 into a cell in Excel -- to force them to stay text instead of being
# auto-converted to numeric values. Column E (Volume(1h)) values already contain
# spaces/% signs so they are never auto-numeric and need no such prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.154.32"
$ws.Range("E2").Value = "  -0.35%  "

# Row 3
$ws.Range("D3").Value = "2.374.00"
$ws.Range("E3").Value = "  -0.89%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'549.51"
$ws.Range("E5").Value = "  +0.04%  "

# Row 6
$ws.Range("D6").Value = "'138.18"
$ws.Range("E6").Value = "  -2.82%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  -1.85%  "

# Row 9
$ws.Range("D9").Value = "2.376.34"
$ws.Range("E9").Value = "  -0.76%  "

# Row 10
$ws.Range("E10").Value = "  +1.35%  "

# Row 11
$ws.Range("E11").Value = "  +1.44%  "

# Row 12
$ws.Range("D12").Value = "'5.37"
$ws.Range("E12").Value = "  +1.29%  "

# Row 13
$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  +0.19%  "

# Row 14
$ws.Range("D14").Value = "'25.06"
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000166"
$ws.Range("E15").Value = "  +0.20%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.784.08"
$ws.Range("E16").Value = "  -1.43%  "

# Row 17
$ws.Range("D17").Value = "61.082.43"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18
$ws.Range("D18").Value = "2.381.10"
$ws.Range("E18").Value = "  -0.54%  "

# Row 19
$ws.Range("D19").Value = "'10.81"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.15"
$ws.Range("E20").Value = "  -0.04%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'320.79"
$ws.Range("E21").Value = "  +0.66%  "

# Row 22
$ws.Range("D22").Value = "'6.69"
$ws.Range("E22").Value = "  -0.53%  "

# Row 23
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("D24").Value = "'64.25"
$ws.Range("E24").Value = "  +0.77%  "

# Row 25
$ws.Range("D25").Value = "'1.68"
$ws.Range("E25").Value = "  -12.88%  "

# Row 26
$ws.Range("D26").Value = "'8.42"
$ws.Range("E26").Value = "  +2.53%  "

# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.12%  "

# Row 28
$ws.Range("D28").Value = "2.477.33"
$ws.Range("E28").Value = "  -1.41%  "

# Row 29
$ws.Range("D29").Value = "'8.16"
$ws.Range("E29").Value = "  +0.76%  "

# Row 30
$ws.Range("D30").Value = "'508.31"
$ws.Range("E30").Value = "  -5.86%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0886"
$ws.Range("E31").Value = "  -6.02%  "

# Row 32
$ws.Range("D32").Value = "'0.150"
$ws.Range("E32").Value = "  +2.63%  "

# Row 33
$ws.Range("D33").Value = "'1.39"
$ws.Range("E33").Value = "  -4.06%  "

# Row 34
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  -0.94%  "

# Row 35
$ws.Range("D35").Value = "'1.52"
$ws.Range("E35").Value = "  -4.10%  "

# Row 36
$ws.Range("E36").Value = "  -0.03%  "

# Row 37
$ws.Range("D37").Value = "'4.71"
$ws.Range("E37").Value = "  -0.39%  "

# Row 38
$ws.Range("D38").Value = "'1.89"
$ws.Range("E38").Value = "  +1.73%  "

# Row 39
$ws.Range("D39").Value = "'0.379"
$ws.Range("E39").Value = "  +0.40%  "

# Row 40
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "'5.36"
$ws.Range("E40").Value = "  -3.52%  "

# Row 41
$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").Value = "'18.59"
$ws.Range("E41").Value = "  +2.52%  "

# Row 42
$ws.Range("D42").Value = "'145.83"
$ws.Range("E42").Value = "  +5.38%  "

# Row 43
$ws.Range("E43").Value = "  -0.08%  "

# Row 44
$ws.Range("D44").Value = "'41.54"
$ws.Range("E44").Value = "  +3.22%  "

# Row 45
$ws.Range("D45").Value = "'148.39"
$ws.Range("E45").Value = "  +4.70%  "

# Row 46
$ws.Range("D46").Value = "'3.61"
$ws.Range("E46").Value = "  -0.57%  "

# Row 47
$ws.Range("D47").Value = "'2.07"
$ws.Range("E47").Value = "  -6.53%  "

# Row 48
$ws.Range("D48").Value = "'0.0521"
$ws.Range("E48").Value = "  +0.02%  "

# Row 49
$ws.Range("D49").Value = "'19.39"
$ws.Range("E49").Value = "  -4.50%  "

# Row 50
$ws.Range("D50").Value = "'0.576"
$ws.Range("E50").Value = "  -0.54%  "

# Row 51
$ws.Range("D51").Value = "'0.0910"
$ws.Range("E51").Value = "  +0.36%  "
